# Add a new "Slovakia" worksheet, cloned from "Portugal", with its own
# market/part-number values, and move the active-sheet/selection state
# from Portugal onto the new Slovakia sheet.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Clone Portugal into a new sheet placed at the end of the tab strip.
$portugal.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# The copied sheet inherited Portugal's custom row heights on rows 3-5;
# reset them back to the sheet's default (auto) height.
$slovakia.Rows("3:5").AutoFit()

# Market-specific values for the new sheet (written in this order so the
# shared-string table gets the part number before the market name).
$slovakia.Range("B4").Value = "NGC-2930/T3236/T3235"
$slovakia.Range("B2").Value = "Slovakia Market"

# Portugal is no longer the active tab; its whole grid is selected instead.
$portugal.Select() | Out-Null
$portugal.Cells.Select() | Out-Null

# Slovakia becomes the active tab, with C13 selected.
$slovakia.Select() | Out-Null
$slovakia.Range("C13").Select() | Out-Null
